$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.479687929153442
$ws.Range("B1").Value = 1.303015947341919
$ws.Range("C1").Value = 4.506392478942871
$ws.Range("D1").Value = 2.119007110595703
$ws.Range("E1").Value = 0.7136600017547607
